$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 381/382 (everything from the old row 381 onward
# shifts down by two, dimension grows from R413 to R415).
$ws.Range("A381:A382").EntireRow.Insert()

# New row 381: "Primera" entry for the new date 2022-08-10 (serial 44783)
$ws.Range("A381").Value = 7
$ws.Range("B381").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C381").Value = "Ñuble"
$ws.Range("D381").Value = 44783
$ws.Range("E381").Value = 16
$ws.Range("F381").Value = 100114014
$ws.Range("G381").Value = "Betarraga"
$ws.Range("H381").Value = "Sin especificar"
$ws.Range("I381").Value = "Primera"
$ws.Range("J381").Value = 200
$ws.Range("K381").Value = 700
$ws.Range("L381").Value = 800
$ws.Range("M381").Value = 750
$ws.Range("N381").Value = "$/paquete 5 unidades"
$ws.Range("O381").Value = "Provincia de Diguillín"
$ws.Range("P381").Value = 150
$ws.Range("Q381").Value = 5
$ws.Range("R381").Value = "Hortaliza"

# New row 382: "Segunda" entry for the same new date
$ws.Range("A382").Value = 7
$ws.Range("B382").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C382").Value = "Ñuble"
$ws.Range("D382").Value = 44783
$ws.Range("E382").Value = 16
$ws.Range("F382").Value = 100114014
$ws.Range("G382").Value = "Betarraga"
$ws.Range("H382").Value = "Sin especificar"
$ws.Range("I382").Value = "Segunda"
$ws.Range("J382").Value = 150
$ws.Range("K382").Value = 600
$ws.Range("L382").Value = 600
$ws.Range("M382").Value = 600
$ws.Range("N382").Value = "$/paquete 5 unidades"
$ws.Range("O382").Value = "Provincia de Diguillín"
$ws.Range("P382").Value = 120
$ws.Range("Q382").Value = 5
$ws.Range("R382").Value = "Hortaliza"
